# Update the "quiz" marksheet: correct marks / total marks values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (B11): correct answer marks per question
$ws.Range("B11").Value = 5

# Total row (B12): total marks scored
$ws.Range("B12").Value = 75

# Total row (E12): "correct/total" summary text
$ws.Range("E12").Value = "75/140"
